$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the label in A6: "microstrain/K" -> "ue/K"
$ws.Range("A6").Value = "Coefficient of thermal expansion (ue/K)"

# Add new columns F (HSMVals) and G (HSMVars) describing where HSM variables
# come from, using the same text style ("@" number format) as columns B:E.
$ws.Range("F1:G6").NumberFormat = "@"

$ws.Range("F1").Value = "HSMVals"
$ws.Range("G1").Value = "HSMVars"

$ws.Range("F2").Value = "1.47"
$ws.Range("G2").Value = "Refractive index"

$ws.Range("F3").Value = "0.527212"
$ws.Range("G3").Value = "Grating period (um)"

$ws.Range("F4").Value = "14"
$ws.Range("G4").Value = "Thermo-optic coeff (ue/K)"

$ws.Range("F5").Value = "20"
$ws.Range("G5").Value = "Ref. temperature (deg C)"

$ws.Range("F6").Value = "0.0015"
$ws.Range("G6").Value = "Fiber length (m)"

$wb.Save()
